# Fixed update to excel issue
# 1) Rename the "Requested quantity" headers on the existing sheets.
# 2) Add a new "PO Forecast" sheet with forecast data (ds, PO_Forecast, yhat_lower, yhat_upper).

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new sheet at the end of the workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Copy header formatting (bold, border, center/top alignment) from an existing header row.
$wsWeekly.Range("A1:B1").Copy($wsForecast.Range("A1:B1"))
$wsWeekly.Range("A1").Copy($wsForecast.Range("C1"))
$wsWeekly.Range("A1").Copy($wsForecast.Range("D1"))

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$data = @(
    @(45333.99999999999, 16, -27.12848536788906, 61.268635485931),
    @(45340.99999999999, 17, -24.89970585588263, 62.86101957473404),
    @(45354.99999999999, 19, -26.0983964814338, 61.98630118097456),
    @(45361.99999999999, 20, -25.28987393518931, 61.6410132474995),
    @(45368.99999999999, 21, -23.13842680542123, 67.10635091518617),
    @(45375.99999999999, 22, -23.48217757148709, 65.46821958965371),
    @(45382.99999999999, 23, -21.2476093949539, 68.19567131024151),
    @(45389.99999999999, 24, -17.29852123269536, 70.23296113937573),
    @(45396.99999999999, 25, -20.60376964184662, 70.58581555609932),
    @(45403.99999999999, 26, -16.74460266930739, 71.82632753207119),
    @(45410.99999999999, 27, -14.32877993647813, 70.77597425706176),
    @(45417.99999999999, 28, -14.98178647506619, 73.21488194057063),
    @(45424.99999999999, 29, -13.90531761314544, 71.94385686426919),
    @(45431.99999999999, 31, -12.83989271506105, 77.3820528008935),
    @(45438.99999999999, 32, -9.57561221294719, 75.94951563564125),
    @(45445.99999999999, 33, -10.44368536844111, 78.45223969845455),
    @(45452.99999999999, 34, -9.692426815856695, 82.83392324319267),
    @(45487.99999999999, 39, -4.935037193975331, 82.19782239226636),
    @(45515.99999999999, 43, -0.01254704164976103, 87.52179427689703),
    @(45522.99999999999, 44, -0.2329252439527258, 86.95515888339557),
    @(45529.99999999999, 45, 3.787039643170443, 89.8833937087519),
    @(45536.99999999999, 46, 3.656492289023299, 90.39751129701759),
    @(45543.99999999999, 47, 5.636410044808433, 93.53662739287346),
    @(45550.99999999999, 48, 5.203478211892017, 93.4191313007082),
    @(45557.99999999999, 49, 4.36945653843658, 92.86456735283494),
    @(45564.99999999999, 50, 5.41886036809309, 91.96687802574476),
    @(45571.99999999999, 51, 6.036771827386374, 96.37047019007309)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row++
}

# Adding a sheet makes it the active tab; restore the original active sheet
# so the workbook-level selection state is left untouched by this edit.
$wsWeekly.Activate()

Write-Output "Sheets now: $($wb.Worksheets.Count)"
